$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily auto-push: insert a new sample row for 2026/02/20 16:00 (ranking 60)
# right before the existing row 821, shifting the remaining rows down by one.
$ws.Rows(821).Insert()

# Force column A to be stored as plain text so the date-like string
# "2026/02/20" isn't auto-converted into a date serial number, then clear
# the formatting change back off so the new cells carry no style index
# (matching the unstyled data rows around them).
$ws.Range("A821").NumberFormat = "@"
$ws.Range("A821").Value = "2026/02/20"
$ws.Range("A821").ClearFormats()
$ws.Range("B821").Value = "金"
$ws.Range("C821").Value = 16
$ws.Range("D821").Value = 60
